# Apply the built-in "Office Theme" design (colour scheme) to the deck's
# slide master, replacing the "Integral" theme that is currently in force.
#
# PowerPoint's RGB values are packed as 0x00BBGGRR, so convert each target
# hex colour (RRGGBB, as used in the OOXML <a:srgbClr val="RRGGBB"/>) into
# that packed integer before assigning it.
function ToBgr($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$p  = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# Office Theme colour scheme (Dark1, Light1, Dark2, Light2, Accent1-6,
# Hyperlink, FollowedHyperlink) - PowerPoint ColorScheme index order.
$cs.Colors(1).RGB  = (ToBgr 0x00 0x00 0x00)   # Dark 1    - 000000
$cs.Colors(2).RGB  = (ToBgr 0xFF 0xFF 0xFF)   # Light 1   - FFFFFF
$cs.Colors(3).RGB  = (ToBgr 0x44 0x54 0x6A)   # Dark 2    - 44546A
$cs.Colors(4).RGB  = (ToBgr 0xE7 0xE6 0xE6)   # Light 2   - E7E6E6
$cs.Colors(5).RGB  = (ToBgr 0x5B 0x9B 0xD5)   # Accent 1  - 5B9BD5
$cs.Colors(6).RGB  = (ToBgr 0xED 0x7D 0x31)   # Accent 2  - ED7D31
$cs.Colors(7).RGB  = (ToBgr 0xA5 0xA5 0xA5)   # Accent 3  - A5A5A5
$cs.Colors(8).RGB  = (ToBgr 0xFF 0xC0 0x00)   # Accent 4  - FFC000
$cs.Colors(9).RGB  = (ToBgr 0x44 0x72 0xC4)   # Accent 5  - 4472C4
$cs.Colors(10).RGB = (ToBgr 0x70 0xAD 0x47)   # Accent 6  - 70AD47
$cs.Colors(11).RGB = (ToBgr 0x05 0x63 0xC1)   # Hyperlink - 0563C1
$cs.Colors(12).RGB = (ToBgr 0x95 0x4F 0x72)   # Followed Hyperlink - 954F72
